$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 190.1
$ws.Range("I4").Value = 156.375
$ws.Range("K4").Value = 156.375
$ws.Range("M4").Value = -42.375

# Row 32
$ws.Range("H32").Value = 568.8570999999999
$ws.Range("I32").Value = 633.3333
$ws.Range("J32").Value = 520.5
$ws.Range("K32").Value = 633.3333
$ws.Range("L32").Value = 520.5
$ws.Range("M32").Value = -307.3333
$ws.Range("N32").Value = -1172.5

# Row 62
$ws.Range("H62").Value = 2049.9048
$ws.Range("I62").Value = 1231.3846
$ws.Range("J62").Value = 3380
$ws.Range("K62").Value = 1231.3846
$ws.Range("L62").Value = 3380
$ws.Range("M62").Value = -607.3846000000001
$ws.Range("N62").Value = -4628

# Row 65
$ws.Range("H65").Value = 2049.9048
$ws.Range("I65").Value = 1231.3846
$ws.Range("J65").Value = 3380
$ws.Range("K65").Value = 6156.923000000001
$ws.Range("L65").Value = 16900
$ws.Range("M65").Value = -3036.923000000001
$ws.Range("N65").Value = -23140

# Row 138
$ws.Range("H138").Value = 3012.4
$ws.Range("J138").Value = 3345.694
$ws.Range("L138").Value = 10037.082
$ws.Range("N138").Value = -20317.082

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 13983.333
$ws.Range("I26").Value = 1975
$ws.Range("K26").Value = 1975
$ws.Range("M26").Value = -1645

# Row 32
$ws.Range("H32").Value = 4007.0657
$ws.Range("I32").Value = 3052.4365
$ws.Range("J32").Value = 17562.8
$ws.Range("K32").Value = 3052.4365
$ws.Range("L32").Value = 17562.8
$ws.Range("M32").Value = -2765.4365
$ws.Range("N32").Value = -18136.8

# Row 34
$ws.Range("H34").Value = 29341.666
$ws.Range("I34").Value = 24012.5
$ws.Range("K34").Value = 24012.5
$ws.Range("M34").Value = -23741.5

# Row 35
$ws.Range("H35").Value = 17340.4
$ws.Range("I35").Value = 6679
$ws.Range("J35").Value = 33332.5
$ws.Range("K35").Value = 6679
$ws.Range("L35").Value = 33332.5
$ws.Range("M35").Value = -6273
$ws.Range("N35").Value = -34144.5

# Row 36
$ws.Range("H36").Value = 35000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 35000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 35000
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -35692

# Row 38
$ws.Range("H38").Value = 8767.6
$ws.Range("I38").Value = 4709.5
$ws.Range("J38").Value = 25000
$ws.Range("K38").Value = 4709.5
$ws.Range("L38").Value = 25000
$ws.Range("M38").Value = -4242.5
$ws.Range("N38").Value = -25934

# Row 39
$ws.Range("H39").Value = 20318.715
$ws.Range("I39").Value = 12446.4
$ws.Range("K39").Value = 12446.4
$ws.Range("M39").Value = -11926.4

# Row 41
$ws.Range("H41").Value = 22234.908
$ws.Range("I41").Value = 3528
$ws.Range("J41").Value = 26392
$ws.Range("K41").Value = 3528
$ws.Range("L41").Value = 26392
$ws.Range("M41").Value = -3114
$ws.Range("N41").Value = -27220

# Row 54
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51538

# Row 74
$ws.Range("H74").Value = 2486.76
$ws.Range("I74").Value = 2573.6924
$ws.Range("K74").Value = 2573.6924
$ws.Range("M74").Value = -1699.6924

# Row 77
$ws.Range("H77").Value = 2486.76
$ws.Range("I77").Value = 2573.6924
$ws.Range("K77").Value = 12868.462
$ws.Range("M77").Value = -8500.462

# Row 132
$ws.Range("H132").Value = 1795.431
$ws.Range("I132").Value = 1193.5581
$ws.Range("J132").Value = 3520.8
$ws.Range("K132").Value = 3580.6743
$ws.Range("L132").Value = 10562.4
$ws.Range("M132").Value = -1050.6743
$ws.Range("N132").Value = -15622.4

$ws = $wb.Worksheets.Item("BSM")
# Row 38
$ws.Range("H38").Value = 21035.5
$ws.Range("J38").Value = 21035.5
$ws.Range("L38").Value = 21035.5
$ws.Range("N38").Value = -21867.5

# Row 44
$ws.Range("H44").Value = 32723.75
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 32723.75
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 32723.75
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -33717.75

# Row 45
$ws.Range("H45").Value = 25000
$ws.Range("J45").Value = 25000
$ws.Range("L45").Value = 25000
$ws.Range("N45").Value = -26616

# Row 49
$ws.Range("H49").Value = 12021.333
$ws.Range("J49").Value = 12021.333
$ws.Range("L49").Value = 12021.333
$ws.Range("N49").Value = -12499.333

# Row 92
$ws.Range("H92").Value = 59401
$ws.Range("J92").Value = 59401
$ws.Range("L92").Value = 59401
$ws.Range("N92").Value = -64393

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 16670060
$ws.Range("I31").Value = 1934.3572
$ws.Range("J31").Value = 31254670
$ws.Range("K31").Value = 1934.3572
$ws.Range("L31").Value = 31254670
$ws.Range("M31").Value = -1639.3572
$ws.Range("N31").Value = -31255260

# Row 34
$ws.Range("H34").Value = 16670060
$ws.Range("I34").Value = 1934.3572
$ws.Range("J34").Value = 31254670
$ws.Range("K34").Value = 1934.3572
$ws.Range("L34").Value = 31254670
$ws.Range("M34").Value = -1732.3572
$ws.Range("N34").Value = -31255074

# Row 134
$ws.Range("H134").Value = 4437.0527
$ws.Range("I134").Value = 5419.727
$ws.Range("J134").Value = 3085.875
$ws.Range("K134").Value = 16259.181
$ws.Range("L134").Value = 9257.625
$ws.Range("M134").Value = -13724.181
$ws.Range("N134").Value = -14327.625

$ws = $wb.Worksheets.Item("CUL")
# Row 41
$ws.Range("H41").Value = 740
$ws.Range("I41").Value = 433.33334
$ws.Range("J41").Value = 1200
$ws.Range("K41").Value = 1300.00002
$ws.Range("L41").Value = 3600
$ws.Range("M41").Value = -962.0000199999999
$ws.Range("N41").Value = -4276

# Row 131
$ws.Range("H131").Value = 5377174.5
$ws.Range("I131").Value = 71429130
$ws.Range("J131").Value = 852.9651
$ws.Range("K131").Value = 214287390
$ws.Range("L131").Value = 2558.8953
$ws.Range("M131").Value = -214282350
$ws.Range("N131").Value = -12638.8953

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1442.3334
$ws.Range("I22").Value = 1122.6666
$ws.Range("K22").Value = 1122.6666
$ws.Range("M22").Value = -827.6666

# Row 27
$ws.Range("H27").Value = 1442.3334
$ws.Range("I27").Value = 1122.6666
$ws.Range("K27").Value = 1122.6666
$ws.Range("M27").Value = -1015.6666

# Row 46
$ws.Range("H46").Value = 2088
$ws.Range("I46").Value = 1704.2858
$ws.Range("J46").Value = 2471.7144
$ws.Range("K46").Value = 1704.2858
$ws.Range("L46").Value = 2471.7144
$ws.Range("M46").Value = -1516.2858
$ws.Range("N46").Value = -2847.7144

# Row 55
$ws.Range("H55").Value = 527.2727
$ws.Range("I55").Value = 443.14285
$ws.Range("J55").Value = 674.5
$ws.Range("K55").Value = 443.14285
$ws.Range("L55").Value = 674.5
$ws.Range("M55").Value = -270.14285
$ws.Range("N55").Value = -1020.5

# Row 132
$ws.Range("H132").Value = 20121.63
$ws.Range("J132").Value = 10945.363
$ws.Range("L132").Value = 32836.089
$ws.Range("N132").Value = -37896.089

$ws = $wb.Worksheets.Item("WVR")
# Row 56
$ws.Range("H56").Value = 18224.75
$ws.Range("I56").Value = 17035
$ws.Range("J56").Value = 18621.334
$ws.Range("K56").Value = 17035
$ws.Range("L56").Value = 18621.334
$ws.Range("M56").Value = -16321
$ws.Range("N56").Value = -20049.334

# Row 132
$ws.Range("H132").Value = 7409546.5
$ws.Range("I132").Value = 1579.375
$ws.Range("J132").Value = 25644542
$ws.Range("K132").Value = 4738.125
$ws.Range("L132").Value = 76933626
$ws.Range("M132").Value = -2208.125
$ws.Range("N132").Value = -76938686
